$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

# Original runs: "Followed" | " " | "by" | " " | "a" | " " | "picture"
# Target runs:   "Followed " | "by " | "a " | "picture"
# Merge each word-run with the single-space run that follows it by
# absorbing the space into the word run's text, then deleting the
# now-redundant standalone space run.

# Merge "Followed" + " " -> "Followed "
$tr.Characters(1, 8).Text = "Followed "
$tr.Characters(10, 1).Text = ""

# Merge "by" + " " -> "by "
$tr.Characters(10, 2).Text = "by "
$tr.Characters(13, 1).Text = ""

# Merge "a" + " " -> "a "
$tr.Characters(13, 1).Text = "a "
$tr.Characters(15, 1).Text = ""
